{"js": "// Apply the resume edits:\n//  1) \"...on platform, 2018\"            -> \"...on platform, 2019\"\n//  2) \"...Classification at Kaggle, 2018\" -> \"...Classification at Kaggle, 2019\"\n//  3) H2O.ai line: \"  Mountain View, California\" -> \"                    Chennai, India\"\n\n// --- 1) \"on platform, 2018\" -> \"on platform, 2019\" -----------------------\nlet platformResults = context.document.body.search(\"on platform, 2018\", { matchCase: true });\nplatformResults.load(\"text\");\nawait context.sync();\nif (platformResults.items.length > 0) {\n  platformResults.items[0].insertText(\"on platform, 2019\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) \"at Kaggle, 2018\" -> \"at Kaggle, 2019\" ----------------------------\nlet kaggleResults = context.document.body.search(\"at Kaggle, 2018\", { matchCase: true });\nkaggleResults.load(\"text\");\nawait context.sync();\nif (kaggleResults.items.length > 0) {\n  kaggleResults.items[0].insertText(\"at Kaggle, 2019\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 3) \"  Mountain View\" -> \"                    Chennai\" ----------------\nlet cityResults = context.document.body.search(\"  Mountain View\", { matchCase: true });\ncityResults.load(\"text\");\nawait context.sync();\nif (cityResults.items.length > 0) {\n  cityResults.items[0].insertText(\"                    Chennai\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 3b) \"California\" -> \"India\" ------------------------------------------\nlet stateResults = context.document.body.search(\"California\", { matchCase: true });\nstateResults.load(\"text\");\nawait context.sync();\nif (stateResults.items.length > 0) {\n  stateResults.items[0].insertText(\"India\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the resume edits:\n#  1) \"...on platform, 2018\"              -> \"...on platform, 2019\"\n#  2) \"...Classification at Kaggle, 2018\" -> \"...Classification at Kaggle, 2019\"\n#  3) H2O.ai line: \"  Mountain View, California\" -> \"                    Chennai, India\"\n\n$d = $word.ActiveDocument\n\n# --- 1) \"on platform, 2018\" -> \"on platform, 2019\" ------------------------\n$find1 = $d.Content.Find\n$find1.Execute(\"on platform, 2018\", $true, $false, $false, $false, $false, $true, 1, $false, \"on platform, 2019\", 2)\n\n# --- 2) \"at Kaggle, 2018\" -> \"at Kaggle, 2019\" -----------------------------\n$find2 = $d.Content.Find\n$find2.Execute(\"at Kaggle, 2018\", $true, $false, $false, $false, $false, $true, 1, $false, \"at Kaggle, 2019\", 2)\n\n# --- 3) \"  Mountain View\" -> \"                    Chennai\" ----------------\n$find3 = $d.Content.Find\n$find3.Execute(\"  Mountain View\", $true, $false, $false, $false, $false, $true, 1, $false, \"                    Chennai\", 2)\n\n# --- 3b) \"California\" -> \"India\" -------------------------------------------\n$find4 = $d.Content.Find\n$find4.Execute(\"California\", $true, $false, $false, $false, $false, $true, 1, $false, \"India\", 2)\n"}
